# Update "想去人数" (want-to-go count) figures that changed between scrapes.
#
# Sheet "展览" (Exhibition):
#   F4 (南宁·草莓动漫节)            1525 -> 1529
#   F5 (南宁·第一届ANE·DACG动漫嘉年华) 705 -> 707
#   F6 (南宁·布谷鸟动漫展4th)          16 -> 17
#
# Sheet "全部类型" (All types) mirrors the same three events at F4/F6/F7
# (row 5 there is an unrelated 演出 entry interleaved by date):
#   F4 (南宁·草莓动漫节)            1525 -> 1529
#   F6 (南宁·第一届ANE·DACG动漫嘉年华) 705 -> 707
#   F7 (南宁·布谷鸟动漫展4th)          16 -> 17

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value = 1529
$wsExhibition.Range("F5").Value = 707
$wsExhibition.Range("F6").Value = 17

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F4").Value = 1529
$wsAllTypes.Range("F6").Value = 707
$wsAllTypes.Range("F7").Value = 17
